# fix bugs : limit bird date
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 values
$ws.Range("B2").Value = 12.5
$ws.Range("C2").Value = 423
$ws.Range("D2").Value = 82.5

# Update row 3 values
$ws.Range("B3").Value = 5
$ws.Range("E3").Value = "Metal"

# Append a new row 47 with data
$ws.Range("A47").Value = 432
$ws.Range("B47").Value = 33
$ws.Range("C47").Value = 33
$ws.Range("D47").Value = 33
$ws.Range("E47").Value = "Plastic"
